$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6058052434456929
$ws1.Range("C2").Value = 0.5613463626492943
$ws1.Range("D2").Value = 0.9681647940074907
$ws1.Range("E2").Value = 0.7106529209621993
$ws1.Range("F2").Value = 0.8456002616944717
$ws1.Range("G2").Value = 0.9419101674724967
$ws1.Range("H2").Value = 0.813400384351022
$ws1.Range("I2").Value = 517
$ws1.Range("J2").Value = 404
$ws1.Range("K2").Value = 130
$ws1.Range("L2").Value = 17

# ---- Sheet: Classification Report ----
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.8843537414965986
$ws2.Range("C2").Value = 0.2434456928838951
$ws2.Range("D2").Value = 0.381791483113069

$ws2.Range("B3").Value = 0.5613463626492943
$ws2.Range("C3").Value = 0.9681647940074907
$ws2.Range("D3").Value = 0.7106529209621993

$ws2.Range("B4").Value = 0.6058052434456929
$ws2.Range("C4").Value = 0.6058052434456929
$ws2.Range("D4").Value = 0.6058052434456929
$ws2.Range("E4").Value = 0.6058052434456929

$ws2.Range("B5").Value = 0.7228500520729464
$ws2.Range("C5").Value = 0.6058052434456929
$ws2.Range("D5").Value = 0.5462222020376342

$ws2.Range("B6").Value = 0.7228500520729465
$ws2.Range("C6").Value = 0.6058052434456929
$ws2.Range("D6").Value = 0.5462222020376342

# ---- Sheet: Confusion Matrix ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 130
$ws3.Range("C2").Value = 404
$ws3.Range("B3").Value = 17
$ws3.Range("C3").Value = 517
